$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.75%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'12.07%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.330"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.65%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08030"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.48%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.599"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.69%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'27.29%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.84%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1270"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.45%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.08%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09559"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.31%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04528"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'9.01%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.19%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001300"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.59%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04213"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.67%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005809"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.18%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.67%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'5.64%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3476"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.87%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.161"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.42%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1394"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.99%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3087"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.87%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001295"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.56%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004332"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-3.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001351"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.71%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003542"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02691"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.93%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05912"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'10.58%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01080"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'93.00%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008026"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.67%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1462"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.07%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007521"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.52%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007926"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.83%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3214"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.44%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007021"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.98%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.73%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05594"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'26.75%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.71%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.73%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.73%"
$ws.Range("E51").Style = "Normal"
